# Generate Report for Handback
# Updates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns for the e6041e91-... row on both the zh-cn and de-de sheets, now
# that a (stale) handback has come in for that file.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30a12f46222588204b54a626d089bd7438ffc7fc/e2e/e6041e91-67a5-4b12-93bc-b88099509b7f.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28207bfdf5c077413188e3a9c1b4d3115f2959a3/e2e/e6041e91-67a5-4b12-93bc-b88099509b7f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30a12f46222588204b54a626d089bd7438ffc7fc/e2e/e6041e91-67a5-4b12-93bc-b88099509b7f.md."

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "e6041e91-67a5-4b12-93bc-b88099509b7f.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackUrl, "", "", "e6041e91-67a5-4b12-93bc-b88099509b7f.md")
$wsZh.Range("J7").Value = "e6041e91-67a5-4b12-93bc-b88099509b7f.8405fc44ad4d60a3f30d7041e7a7fd8cc89195a4.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-28 18:53:49"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "e6041e91-67a5-4b12-93bc-b88099509b7f.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackUrl, "", "", "e6041e91-67a5-4b12-93bc-b88099509b7f.md")
$wsDe.Range("J7").Value = "e6041e91-67a5-4b12-93bc-b88099509b7f.8405fc44ad4d60a3f30d7041e7a7fd8cc89195a4.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-28 18:53:56"
$wsDe.Range("P7").Value = $errorDetail
